$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 1
    7  = 1
    8  = 1
    9  = 0
    10 = 1
    12 = 1
    13 = 2
    14 = 2
    15 = 1
    16 = 0
    17 = 1
    18 = 2
    19 = 1
    20 = 0
    21 = 1
    23 = 0
    25 = 0
    26 = 2
    27 = 1
    28 = 1
    29 = 1
    31 = 1
    32 = 0
    34 = 1
    35 = 1
    36 = 1
    37 = 1
    38 = 0
    39 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
